# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (want-to-go count) figures pulled from bilibili
# for the convention rows that changed in this scrape.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("F15").Value = 2243
$wsExhibit.Range("F16").Value = 693
$wsExhibit.Range("F17").Value = 13939
$wsExhibit.Range("F19").Value = 1303

# Sheet 2: 演出 (Performances)
$wsShow = $wb.Worksheets.Item(2)
$wsShow.Range("F10").Value = 19
$wsShow.Range("F12").Value = 65

# Sheet 3: 本地生活 (Local Life)
$wsLocal = $wb.Worksheets.Item(3)
$wsLocal.Range("F2").Value = 5730
$wsLocal.Range("F3").Value = 491

# Sheet 4: 全部类型 (All Types)
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F3").Value = 491
$wsAll.Range("F22").Value = 19
$wsAll.Range("F23").Value = 2243
$wsAll.Range("F24").Value = 693
$wsAll.Range("F26").Value = 65
$wsAll.Range("F27").Value = 1303
